# Apply scraped crypto price/volume updates from cryptos.xlsx commit
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'" + '42.731.78'
$ws.Range("E2").Value = '  -0.46%  '
$ws.Range("D3").Value = "'" + '2.520.58'
$ws.Range("E4").Value = '  +0.10%  '
$ws.Range("D5").Value = "'" + '304.15'
$ws.Range("E5").Value = '  +0.45%  '
$ws.Range("D6").Value = "'" + '97.21'
$ws.Range("E6").Value = '  +0.70%  '
$ws.Range("D7").Value = "'" + '0.576'
$ws.Range("E7").Value = '  +0.13%  '
$ws.Range("D9").Value = "'" + '0.540'
$ws.Range("E9").Value = '  -1.62%  '
$ws.Range("D10").Value = "'" + '36.73'
$ws.Range("E10").Value = '  +0.08%  '
$ws.Range("D11").Value = "'" + '0.0811'
$ws.Range("E11").Value = '  -0.28%  '
$ws.Range("E12").Value = '  +0.38%  '
$ws.Range("E13").Value = '  -0.87%  '
$ws.Range("D14").Value = "'" + '2.906.37'
$ws.Range("E14").Value = '  -2.18%  '
$ws.Range("D15").Value = "'" + '2.515.59'
$ws.Range("E15").Value = '  -4.65%  '
$ws.Range("D16").Value = "'" + '15.02'
$ws.Range("E16").Value = '  +5.16%  '
$ws.Range("D17").Value = "'" + '0.860'
$ws.Range("E17").Value = '  -2.92%  '
$ws.Range("D18").Value = "'" + '42.712.26'
$ws.Range("E18").Value = '  -0.69%  '
$ws.Range("D19").Value = "'" + '12.93'
$ws.Range("E19").Value = '  -0.60%  '
$ws.Range("D20").Value = "'" + '0.0₃0973'
$ws.Range("E20").Value = '  -2.35%  '
$ws.Range("D21").Value = "'" + '6.43'
$ws.Range("E21").Value = '  -3.60%  '
$ws.Range("D22").Value = "'" + '71.17'
$ws.Range("D23").Value = "'" + '250.74'
$ws.Range("E23").Value = '  -1.37%  '
$ws.Range("D24").Value = "'" + '2.91'
$ws.Range("E24").Value = '  -1.71%  '
$ws.Range("D25").Value = "'" + '2.03'
$ws.Range("E25").Value = '  -5.25%  '
$ws.Range("D26").Value = "'" + '26.88'
$ws.Range("E26").Value = '  -7.07%  '
$ws.Range("E27").Value = '  +0.01%  '
$ws.Range("E28").Value = '  +11.06%  '
$ws.Range("D29").Value = "'" + '10.33'
$ws.Range("E29").Value = '  +0.01%  '
$ws.Range("D30").Value = "'" + '37.94'
$ws.Range("E30").Value = '  +1.13%  '
$ws.Range("E31").Value = '  -1.03%  '
$ws.Range("D32").Value = "'" + '156.38'
$ws.Range("E32").Value = '  +0.88%  '
$ws.Range("D33").Value = "'" + '0.0790'
$ws.Range("E33").Value = '  -1.79%  '
$ws.Range("E34").Value = '  -4.53%  '
$ws.Range("D35").Value = "'" + '3.26'
$ws.Range("E35").Value = '  -4.75%  '
$ws.Range("E36").Value = '  -4.58%  '
$ws.Range("B37").Value = 'Kaspa'
$ws.Range("C37").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D37").Value = "'" + '0.116'
$ws.Range("E37").Value = '  +1.65%  '
$ws.Range("B38").Value = 'Celestia'
$ws.Range("C38").Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range("D38").Value = "'" + '18.27'
$ws.Range("E38").Value = '  +0.34%  '
$ws.Range("B39").Value = 'Stellar'
$ws.Range("C39").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D39").Value = "'" + '0.119'
$ws.Range("E39").Value = '  -0.90%  '
$ws.Range("B40").Value = 'EnergySwap'
$ws.Range("C40").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D40").Value = "'" + '24.01'
$ws.Range("E40").Value = '  +3.71%  '
$ws.Range("E41").Value = '  -4.07%  '
$ws.Range("D42").Value = "'" + '3.38'
$ws.Range("E42").Value = '  -1.76%  '
$ws.Range("D43").Value = "'" + '3.84'
$ws.Range("E43").Value = '  -1.52%  '
$ws.Range("D44").Value = "'" + '0.999'
$ws.Range("E44").Value = '  +0.01%  '
$ws.Range("D45").Value = "'" + '0.0300'
$ws.Range("E45").Value = '  -3.77%  '
$ws.Range("D46").Value = "'" + '2.028.59'
$ws.Range("E46").Value = '  -2.72%  '
$ws.Range("D47").Value = "'" + '85.26'
$ws.Range("E47").Value = '  -0.09%  '
$ws.Range("D48").Value = "'" + '8.98'
$ws.Range("E48").Value = '  -3.06%  '
$ws.Range("D49").Value = "'" + '2.766.82'
$ws.Range("E49").Value = '  -2.06%  '
$ws.Range("E50").Value = '  -1.37%  '
$ws.Range("D51").Value = "'" + '101.73'
$ws.Range("E51").Value = '  -4.69%  '

# Clear the quote-prefix style the apostrophe trick added to column D cells
$ws.Range("D2").Style = "Normal"
$ws.Range("D3").Style = "Normal"
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").Style = "Normal"
$ws.Range("D7").Style = "Normal"
$ws.Range("D9").Style = "Normal"
$ws.Range("D10").Style = "Normal"
$ws.Range("D11").Style = "Normal"
$ws.Range("D14").Style = "Normal"
$ws.Range("D15").Style = "Normal"
$ws.Range("D16").Style = "Normal"
$ws.Range("D17").Style = "Normal"
$ws.Range("D18").Style = "Normal"
$ws.Range("D19").Style = "Normal"
$ws.Range("D20").Style = "Normal"
$ws.Range("D21").Style = "Normal"
$ws.Range("D22").Style = "Normal"
$ws.Range("D23").Style = "Normal"
$ws.Range("D24").Style = "Normal"
$ws.Range("D25").Style = "Normal"
$ws.Range("D26").Style = "Normal"
$ws.Range("D29").Style = "Normal"
$ws.Range("D30").Style = "Normal"
$ws.Range("D32").Style = "Normal"
$ws.Range("D33").Style = "Normal"
$ws.Range("D35").Style = "Normal"
$ws.Range("D37").Style = "Normal"
$ws.Range("D38").Style = "Normal"
$ws.Range("D39").Style = "Normal"
$ws.Range("D40").Style = "Normal"
$ws.Range("D42").Style = "Normal"
$ws.Range("D43").Style = "Normal"
$ws.Range("D44").Style = "Normal"
$ws.Range("D45").Style = "Normal"
$ws.Range("D46").Style = "Normal"
$ws.Range("D47").Style = "Normal"
$ws.Range("D48").Style = "Normal"
$ws.Range("D49").Style = "Normal"
$ws.Range("D51").Style = "Normal"
